$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enoch (row 7) paid for May,18 (E7: 300 -> 400) and also paid June,18 (F7: empty -> 100)
$ws.Range("E7").Value = 400
$ws.Range("F7").Value = 100

# Booked slot #8 for 31/05/2018, 6-7 in the Booking History table (row 47)
$ws.Range("B47").Value = "31/05/2018, 6-7"

# Booking date (29/05/2018, serial 43249) - copy the date formatting from the row above
$ws.Range("C46").Copy()
$ws.Range("C47").PasteSpecial(-4122)
$ws.Range("C47").Value = 43249
$excel.CutCopyMode = $false

# Amount paid for the slot
$ws.Range("D47").Value = 1265

# Move the view roughly where the author left it, and select D47
$ws.Range("D47").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
